{"js": "// 1. Update the letter date from September 19, 2025 to September 21, 2025.\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n}\n\n// 2. Split the single-line mailing address into two lines:\n//    \"2900 Sanor Pl, Santa Clara CA 95051\" ->\n//      \"2900 Sanor Pl\"\n//      \"Santa Clara, CA 95051\"\nconst addressResults = context.document.body.search(\"2900 Sanor Pl, Santa Clara CA 95051\", { matchCase: true });\naddressResults.load(\"text\");\nawait context.sync();\n\nif (addressResults.items.length > 0) {\n  const addressRange = addressResults.items[0];\n  const addressParagraph = addressRange.paragraphs.getFirst();\n  addressParagraph.load(\"text\");\n  await context.sync();\n\n  // Insert the new second line of the address right after the existing paragraph,\n  // inheriting the same paragraph/run formatting.\n  addressParagraph.insertParagraph(\"Santa Clara, CA 95051\", Word.InsertLocation.after);\n\n  // Trim the original paragraph down to just the street address.\n  addressRange.insertText(\"2900 Sanor Pl\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3. Remove the now-superfluous empty \"No Spacing\" paragraph that directly\n//    follows the \"Board of Directors\" signature line.\nconst boardResults = context.document.body.search(\"Board of Directors\", { matchCase: true });\nboardResults.load(\"text\");\nawait context.sync();\n\nif (boardResults.items.length > 0) {\n  const boardParagraph = boardResults.items[0].paragraphs.getFirst();\n  boardParagraph.load(\"text\");\n  await context.sync();\n\n  const nextParagraph = boardParagraph.getNext();\n  nextParagraph.load(\"text\");\n  await context.sync();\n\n  if (nextParagraph.text === \"\") {\n    nextParagraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Update the letter date from September 19, 2025 to September 21, 2025.\n$find = $d.Content.Find\n$found = $find.Execute(\"September 19, 2025\", $false, $false, $false, $false, $false, $true, 1, $false, \"September 21, 2025\", 2)\n\n# 2. Split the single-line mailing address into two lines:\n#    \"2900 Sanor Pl, Santa Clara CA 95051\" ->\n#      \"2900 Sanor Pl\"\n#      \"Santa Clara, CA 95051\"\n$addrRange = $d.Content\n$addrFind = $addrRange.Find\n$addrFound = $addrFind.Execute(\"2900 Sanor Pl, Santa Clara CA 95051\")\nif ($addrFound) {\n    $addrRange.Text = \"2900 Sanor Pl\" + [char]13 + \"Santa Clara, CA 95051\"\n}\n\n# 3. Remove the now-superfluous empty \"No Spacing\" paragraph that directly\n#    follows the \"Board of Directors\" signature line.\n$boardRange = $d.Content\n$boardFind = $boardRange.Find\n$boardFound = $boardFind.Execute(\"Board of Directors\")\nif ($boardFound) {\n    $boardPara = $boardRange.Paragraphs(1)\n    $nextPara = $boardPara.Next()\n    if ($nextPara.Range.Text.Trim() -eq \"\") {\n        $nextPara.Range.Delete()\n    }\n}\n"}
